# Natmi following Dr Hou advice
# Update Arf1-Chrm3 LR-pairs data: recompute existing rows 2-7 with new
# ligand/receptor-expressing-cell counts (grouped by 3 samples) and add
# 3 new rows (8-10) for the sCs -> {ECs,FAPs,sCs} target-cluster combos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Arf1"
$ws.Cells.Item(2, 3).Value = "Chrm3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 77.63718533333333
$ws.Cells.Item(2, 8).Value = 232.911556
$ws.Cells.Item(2, 9).Value = 0.4513549673384918
$ws.Cells.Item(2, 10).Value = 0.4513549673384918
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 2.425673333333334
$ws.Cells.Item(2, 14).Value = 7.27702
$ws.Cells.Item(2, 15).Value = 0.6040766302760682
$ws.Cells.Item(2, 16).Value = 0.6040766302760682
$ws.Cells.Item(2, 17).Value = 188.3224501381245
$ws.Cells.Item(2, 18).Value = 1694.90205124312
$ws.Cells.Item(2, 19).Value = 0.272652987728201
$ws.Cells.Item(2, 20).Value = 0.2726529877282009

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Arf1"
$ws.Cells.Item(3, 3).Value = "Chrm3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 77.63718533333333
$ws.Cells.Item(3, 8).Value = 232.911556
$ws.Cells.Item(3, 9).Value = 0.4513549673384918
$ws.Cells.Item(3, 10).Value = 0.4513549673384918
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.9912810000000002
$ws.Cells.Item(3, 14).Value = 2.973843
$ws.Cells.Item(3, 15).Value = 0.2468632844777222
$ws.Cells.Item(3, 16).Value = 0.2468632844777221
$ws.Cells.Item(3, 17).Value = 76.96026671441201
$ws.Cells.Item(3, 18).Value = 692.6424004297081
$ws.Cells.Item(3, 19).Value = 0.1114229697025151
$ws.Cells.Item(3, 20).Value = 0.1114229697025151

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Arf1"
$ws.Cells.Item(4, 3).Value = "Chrm3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 77.63718533333333
$ws.Cells.Item(4, 8).Value = 232.911556
$ws.Cells.Item(4, 9).Value = 0.4513549673384918
$ws.Cells.Item(4, 10).Value = 0.4513549673384918
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.5985516666666667
$ws.Cells.Item(4, 14).Value = 1.795655
$ws.Cells.Item(4, 15).Value = 0.1490600852462097
$ws.Cells.Item(4, 16).Value = 0.1490600852462097
$ws.Cells.Item(4, 17).Value = 46.46986667657556
$ws.Cells.Item(4, 18).Value = 418.22880008918
$ws.Cells.Item(4, 19).Value = 0.0672790099077758
$ws.Cells.Item(4, 20).Value = 0.0672790099077758

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Arf1"
$ws.Cells.Item(5, 3).Value = "Chrm3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 62.40815866666667
$ws.Cells.Item(5, 8).Value = 187.224476
$ws.Cells.Item(5, 9).Value = 0.3628188257432201
$ws.Cells.Item(5, 10).Value = 0.3628188257432201
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 2.425673333333334
$ws.Cells.Item(5, 14).Value = 7.27702
$ws.Cells.Item(5, 15).Value = 0.6040766302760682
$ws.Cells.Item(5, 16).Value = 0.6040766302760682
$ws.Cells.Item(5, 17).Value = 151.3818062601689
$ws.Cells.Item(5, 18).Value = 1362.43625634152
$ws.Cells.Item(5, 19).Value = 0.2191703736556844
$ws.Cells.Item(5, 20).Value = 0.2191703736556844

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Arf1"
$ws.Cells.Item(6, 3).Value = "Chrm3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 62.40815866666667
$ws.Cells.Item(6, 8).Value = 187.224476
$ws.Cells.Item(6, 9).Value = 0.3628188257432201
$ws.Cells.Item(6, 10).Value = 0.3628188257432201
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9912810000000002
$ws.Cells.Item(6, 14).Value = 2.973843
$ws.Cells.Item(6, 15).Value = 0.2468632844777222
$ws.Cells.Item(6, 16).Value = 0.2468632844777221
$ws.Cells.Item(6, 17).Value = 61.86402193125202
$ws.Cells.Item(6, 18).Value = 556.7761973812682
$ws.Cells.Item(6, 19).Value = 0.08956664699332165
$ws.Cells.Item(6, 20).Value = 0.08956664699332163

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Arf1"
$ws.Cells.Item(7, 3).Value = "Chrm3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 62.40815866666667
$ws.Cells.Item(7, 8).Value = 187.224476
$ws.Cells.Item(7, 9).Value = 0.3628188257432201
$ws.Cells.Item(7, 10).Value = 0.3628188257432201
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.5985516666666667
$ws.Cells.Item(7, 14).Value = 1.795655
$ws.Cells.Item(7, 15).Value = 0.1490600852462097
$ws.Cells.Item(7, 16).Value = 0.1490600852462097
$ws.Cells.Item(7, 17).Value = 37.35450738353111
$ws.Cells.Item(7, 18).Value = 336.19056645178
$ws.Cells.Item(7, 19).Value = 0.0540818050942141
$ws.Cells.Item(7, 20).Value = 0.0540818050942141

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Arf1"
$ws.Cells.Item(8, 3).Value = "Chrm3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.96380833333333
$ws.Cells.Item(8, 8).Value = 95.891425
$ws.Cells.Item(8, 9).Value = 0.1858262069182881
$ws.Cells.Item(8, 10).Value = 0.1858262069182881
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 2.425673333333334
$ws.Cells.Item(8, 14).Value = 7.27702
$ws.Cells.Item(8, 15).Value = 0.6040766302760682
$ws.Cells.Item(8, 16).Value = 0.6040766302760682
$ws.Cells.Item(8, 17).Value = 77.53375750594445
$ws.Cells.Item(8, 18).Value = 697.8038175535
$ws.Cells.Item(8, 19).Value = 0.1122532688921828
$ws.Cells.Item(8, 20).Value = 0.1122532688921828

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Arf1"
$ws.Cells.Item(9, 3).Value = "Chrm3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.96380833333333
$ws.Cells.Item(9, 8).Value = 95.891425
$ws.Cells.Item(9, 9).Value = 0.1858262069182881
$ws.Cells.Item(9, 10).Value = 0.1858262069182881
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.9912810000000002
$ws.Cells.Item(9, 14).Value = 2.973843
$ws.Cells.Item(9, 15).Value = 0.2468632844777222
$ws.Cells.Item(9, 16).Value = 0.2468632844777221
$ws.Cells.Item(9, 17).Value = 31.68511588847501
$ws.Cells.Item(9, 18).Value = 285.1660429962751
$ws.Cells.Item(9, 19).Value = 0.04587366778188541
$ws.Cells.Item(9, 20).Value = 0.0458736677818854

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Arf1"
$ws.Cells.Item(10, 3).Value = "Chrm3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 31.96380833333333
$ws.Cells.Item(10, 8).Value = 95.891425
$ws.Cells.Item(10, 9).Value = 0.1858262069182881
$ws.Cells.Item(10, 10).Value = 0.1858262069182881
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.5985516666666667
$ws.Cells.Item(10, 14).Value = 1.795655
$ws.Cells.Item(10, 15).Value = 0.1490600852462097
$ws.Cells.Item(10, 16).Value = 0.1490600852462097
$ws.Cells.Item(10, 17).Value = 19.13199075093056
$ws.Cells.Item(10, 18).Value = 172.187916758375
$ws.Cells.Item(10, 19).Value = 0.02769927024421983
$ws.Cells.Item(10, 20).Value = 0.02769927024421983
